$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated price / volume values scraped on 2023-02-11 13:39 UTC
# Ensure cells remain text (matching source data format: inline strings)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "309.27"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.85%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.21"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "5.02%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.134"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.33%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07629"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.55%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.39%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "2.477"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.76%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9063"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.93%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1113"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "8.01%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1799"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.28%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09085"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.08%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04261"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.91%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1050"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.36%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001259"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.88%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005689"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.50%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.342"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.46%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.250"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.35%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.18%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.689"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-4.85%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.25%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2707"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.89%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04030"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-2.74%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001258"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "4.33%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004094"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.72%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001299"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.16%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003743"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02420"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "1.69%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05249"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1.72%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007787"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.46%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1302"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.01%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007047"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "20.66%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.12%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008450"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "14.08%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3337"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.32%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006889"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "6.93%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05481"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "1,214.88%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "39.85%"
